# Auto-generated Excel COM-interop script to apply the Gungnir_Profits.xlsx diff
# Updates currentAveragePrice / Leve price / profit columns (H-N) on multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 10001
$ws.Range("I70").Value = 10001
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 30003
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -29733
$ws.Range("H73").Value = 10001
$ws.Range("I73").Value = 10001
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 30003
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -29067
$ws.Range("H92").Value = 1401.6666
$ws.Range("I92").Value = 2360
$ws.Range("J92").Value = 443.33334
$ws.Range("K92").Value = 2360
$ws.Range("L92").Value = 443.33334
$ws.Range("M92").Value = -1112
$ws.Range("N92").Value = -2939.33334
$ws.Range("H98").Value = 156250640
$ws.Range("I98").Value = 178571870
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 178571870
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -178570372
$ws.Range("N98").Value = -4996
$ws.Range("H112").Value = 16404.25
$ws.Range("I112").Value = 475
$ws.Range("J112").Value = 17041.42
$ws.Range("K112").Value = 1425
$ws.Range("L112").Value = 51124.25999999999
$ws.Range("M112").Value = -317
$ws.Range("N112").Value = -53340.25999999999
$ws.Range("H122").Value = 156250640
$ws.Range("I122").Value = 178571870
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 535715610
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -535713160
$ws.Range("N122").Value = -10900
$ws.Range("H137").Value = 1125.0167
$ws.Range("I137").Value = 977.2553
$ws.Range("J137").Value = 1659.2307
$ws.Range("K137").Value = 2931.7659
$ws.Range("L137").Value = 4977.6921
$ws.Range("M137").Value = -381.7659000000003
$ws.Range("N137").Value = -10077.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1094.1923
$ws.Range("I61").Value = 1064.6285
$ws.Range("J61").Value = 1155.0588
$ws.Range("K61").Value = 1064.6285
$ws.Range("L61").Value = 1155.0588
$ws.Range("M61").Value = -852.6285
$ws.Range("N61").Value = -1579.0588
$ws.Range("H74").Value = 1242.1052
$ws.Range("I74").Value = 1386
$ws.Range("J74").Value = 888.9091
$ws.Range("K74").Value = 1386
$ws.Range("L74").Value = 888.9091
$ws.Range("M74").Value = -512
$ws.Range("N74").Value = -2636.9091
$ws.Range("H77").Value = 1242.1052
$ws.Range("I77").Value = 1386
$ws.Range("J77").Value = 888.9091
$ws.Range("K77").Value = 6930
$ws.Range("L77").Value = 4444.5455
$ws.Range("M77").Value = -2562
$ws.Range("N77").Value = -13180.5455
$ws.Range("H111").Value = 17950
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 17950
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 17950
$ws.Range("N111").Value = -26130
$ws.Range("H136").Value = 1094.1923
$ws.Range("I136").Value = 1064.6285
$ws.Range("J136").Value = 1155.0588
$ws.Range("K136").Value = 3193.8855
$ws.Range("L136").Value = 3465.1764
$ws.Range("M136").Value = -643.8855000000003
$ws.Range("N136").Value = -8565.1764

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3271872.8
$ws.Range("I134").Value = 826.2143
$ws.Range("J134").Value = 18536756
$ws.Range("K134").Value = 2478.6429
$ws.Range("L134").Value = 55610268
$ws.Range("M134").Value = 56.35710000000017
$ws.Range("N134").Value = -55615338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1241.5
$ws.Range("I31").Value = 951.1177
$ws.Range("J31").Value = 1790
$ws.Range("K31").Value = 951.1177
$ws.Range("L31").Value = 1790
$ws.Range("M31").Value = -656.1177
$ws.Range("N31").Value = -2380
$ws.Range("H34").Value = 1241.5
$ws.Range("I34").Value = 951.1177
$ws.Range("J34").Value = 1790
$ws.Range("K34").Value = 951.1177
$ws.Range("L34").Value = 1790
$ws.Range("M34").Value = -749.1177
$ws.Range("N34").Value = -2194
$ws.Range("H58").Value = 17858398
$ws.Range("I58").Value = 22223584
$ws.Range("J58").Value = 827.2727
$ws.Range("K58").Value = 22223584
$ws.Range("L58").Value = 827.2727
$ws.Range("M58").Value = -22223381
$ws.Range("N58").Value = -1233.2727
$ws.Range("H132").Value = 33335374
$ws.Range("I132").Value = 1950
$ws.Range("J132").Value = 83335510
$ws.Range("K132").Value = 5850
$ws.Range("L132").Value = 250006530
$ws.Range("M132").Value = -3320
$ws.Range("N132").Value = -250011590
$ws.Range("H134").Value = 1175.1666
$ws.Range("I134").Value = 1182.3572
$ws.Range("J134").Value = 1150
$ws.Range("K134").Value = 3547.0716
$ws.Range("L134").Value = 3450
$ws.Range("M134").Value = -1012.0716
$ws.Range("N134").Value = -8520
$ws.Range("H136").Value = 17858398
$ws.Range("I136").Value = 22223584
$ws.Range("J136").Value = 827.2727
$ws.Range("K136").Value = 66670752
$ws.Range("L136").Value = 2481.8181
$ws.Range("M136").Value = -66668202
$ws.Range("N136").Value = -7581.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3590.875
$ws.Range("I3").Value = 3590.875
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10772.625
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -10660.625
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 40234120
$ws.Range("I5").Value = 60606668
$ws.Range("J5").Value = 27784226
$ws.Range("K5").Value = 181820004
$ws.Range("L5").Value = 83352678
$ws.Range("M5").Value = -181819892
$ws.Range("N5").Value = -83352902
$ws.Range("H122").Value = 8776337
$ws.Range("I122").Value = 41666960
$ws.Range("J122").Value = 5503.8887
$ws.Range("K122").Value = 375002640
$ws.Range("L122").Value = 49534.99830000001
$ws.Range("M122").Value = -375000190
$ws.Range("H129").Value = 1145.625
$ws.Range("I129").Value = 821.0714
$ws.Range("J129").Value = 1600
$ws.Range("K129").Value = 2463.2142
$ws.Range("L129").Value = 4800
$ws.Range("M129").Value = 2536.7858
$ws.Range("N129").Value = -14800
$ws.Range("H130").Value = 2266.2
$ws.Range("I130").Value = 530
$ws.Range("J130").Value = 2459.111
$ws.Range("K130").Value = 1590
$ws.Range("L130").Value = 7377.333
$ws.Range("M130").Value = 3430
$ws.Range("N130").Value = -17417.333
$ws.Range("H131").Value = 818
$ws.Range("I131").Value = 529.75
$ws.Range("J131").Value = 843.06525
$ws.Range("K131").Value = 1589.25
$ws.Range("L131").Value = 2529.19575
$ws.Range("M131").Value = 3450.75
$ws.Range("N131").Value = -12609.19575
$ws.Range("H135").Value = 40234120
$ws.Range("I135").Value = 60606668
$ws.Range("J135").Value = 27784226
$ws.Range("K135").Value = 545460012
$ws.Range("L135").Value = 250058034
$ws.Range("M135").Value = -545457477
$ws.Range("N135").Value = -250063104
$ws.Range("H136").Value = 35716940
$ws.Range("I136").Value = 56819100
$ws.Range("J136").Value = 5599.154
$ws.Range("K136").Value = 170457300
$ws.Range("L136").Value = 16797.462
$ws.Range("M136").Value = -170452200
$ws.Range("N136").Value = -26997.462
$ws.Range("H137").Value = 46297330
$ws.Range("I137").Value = 38462556
$ws.Range("J137").Value = 66667736
$ws.Range("K137").Value = 115387668
$ws.Range("L137").Value = 200003208
$ws.Range("M137").Value = -115382568
$ws.Range("N137").Value = -200013408
$ws.Range("H139").Value = 11606051
$ws.Range("I139").Value = 16667548
$ws.Range("J139").Value = 759985.7
$ws.Range("K139").Value = 50002644
$ws.Range("L139").Value = 2279957.1
$ws.Range("M139").Value = -49997504
$ws.Range("N139").Value = -2290237.1
$ws.Range("H140").Value = 30002310
$ws.Range("I140").Value = 34617280
$ws.Range("J140").Value = 4998.5
$ws.Range("K140").Value = 103851840
$ws.Range("L140").Value = 14995.5
$ws.Range("M140").Value = -103846660
$ws.Range("N140").Value = -25355.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 5000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5384
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 14994.214
$ws.Range("I132").Value = 8455.308000000001
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 25365.924
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -22835.924
$ws.Range("N132").Value = -305060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18523998
$ws.Range("I132").Value = 34484644
$ws.Range("J132").Value = 9652.200000000001
$ws.Range("K132").Value = 103453932
$ws.Range("L132").Value = 28956.6
$ws.Range("M132").Value = -103451402
$ws.Range("N132").Value = -34016.60000000001
$ws.Range("H136").Value = 1746.9111
$ws.Range("I136").Value = 1966.5161
$ws.Range("J136").Value = 1260.6428
$ws.Range("K136").Value = 5899.5483
$ws.Range("L136").Value = 3781.9284
$ws.Range("M136").Value = -3349.5483
$ws.Range("N136").Value = -8881.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 18000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 18000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 18000
$ws.Range("N64").Value = -18496
$ws.Range("H67").Value = 18000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 18000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 18000
$ws.Range("N67").Value = -19716
$ws.Range("H75").Value = 11266.667
$ws.Range("I75").Value = 11500
$ws.Range("J75").Value = 10800
$ws.Range("K75").Value = 11500
$ws.Range("L75").Value = 10800
$ws.Range("M75").Value = -10564
$ws.Range("N75").Value = -12672
$ws.Range("H78").Value = 11266.667
$ws.Range("I78").Value = 11500
$ws.Range("J78").Value = 10800
$ws.Range("K78").Value = 34500
$ws.Range("L78").Value = 32400
$ws.Range("M78").Value = -29820
$ws.Range("N78").Value = -41760
$ws.Range("H132").Value = 21333.719
$ws.Range("I132").Value = 30455.914
$ws.Range("J132").Value = 6821.136
$ws.Range("K132").Value = 91367.742
$ws.Range("L132").Value = 20463.408
$ws.Range("M132").Value = -88837.742
$ws.Range("N132").Value = -25523.408
$ws.Range("H136").Value = 9263748
$ws.Range("I136").Value = 15157702
$ws.Range("J136").Value = 1819.1904
$ws.Range("K136").Value = 45473106
$ws.Range("L136").Value = 5457.5712
$ws.Range("M136").Value = -45470556
$ws.Range("N136").Value = -10557.5712
